# This workbook's data rows (3-7) got rotated: each row's full content moved
# down by one row, with row 7's content wrapping around to row 3.
#   new row 3 <- old row 7
#   new row 4 <- old row 3
#   new row 5 <- old row 4
#   new row 6 <- old row 5
#   new row 7 <- old row 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol = "AY"

# All columns in use on the sheet, in order, matching the Value2 array layout
# (1-based index into the array below corresponds to this list's position).
$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY")

# Capture the full content (as variant arrays) of every source row before
# writing anything back, so that overwriting a row doesn't clobber data that
# still needs to be read for another row.
$data = @{}
foreach ($r in 3..7) {
    $data[$r] = $ws.Range("$firstCol$r`:$lastCol$r").Value2
}

# Mapping of destination row -> source row, implementing the rotation.
$mapping = @{ 3 = 7; 4 = 3; 5 = 4; 6 = 5; 7 = 6 }

# Columns that may hold date/time strings formatted like "2014-08-14" or
# "00:00"; these must stay plain text instead of being auto-converted to
# Excel date/time serial numbers when assigned.
$dateLikeCols = @{ "Y" = $true; "Z" = $true; "AA" = $true; "AB" = $true }

# Column I ("Antal") holds a numeral stored as text ("1") on one of the
# rows; Value2 can't distinguish that from a real number, so it must be
# written back as an explicit string.
$textNumberCols = @{ "I" = $true }

foreach ($destRow in 3..7) {
    $srcRow = $mapping[$destRow]
    $values = $data[$srcRow]

    # Clear the destination row completely first so that columns which exist
    # in the destination but not in the source (e.g. Z/AB on the row that
    # used to be the "no time" row) don't keep stale values.
    $ws.Range("$firstCol$destRow`:$lastCol$destRow").ClearContents()

    for ($i = 0; $i -lt $allCols.Length; $i++) {
        $col = $allCols[$i]
        $val = $values[1, $i + 1]

        if ($val -eq $null) {
            continue
        }
        if (($val -is [string]) -and ($val -eq "")) {
            continue
        }

        $cell = $ws.Range("$col$destRow")

        if ($dateLikeCols.ContainsKey($col) -or $textNumberCols.ContainsKey($col)) {
            $cell.NumberFormat = "@"
            $cell.Value = [string]$val
        }
        else {
            $cell.Value = $val
        }
    }
}
